$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.743794441223145
$ws.Range("B1").Value = 4.221025466918945
$ws.Range("C1").Value = 3.159626483917236
$ws.Range("D1").Value = 2.159981727600098
$ws.Range("E1").Value = 1.946651220321655
